$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3-6 get their D,K,L,M,N,O,P,R,S values cyclically re-assigned:
#   new row3 <- old row5, new row4 <- old row6, new row5 <- old row3, new row6 <- old row4

$ws.Range("D3").Value = 44305
$ws.Range("K3").Value = "Mankaki"
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 250
$ws.Range("N3").Value = 24000
$ws.Range("O3").Value = 25000
$ws.Range("P3").Value = 24500
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 1361

$ws.Range("D4").Value = 44313
$ws.Range("K4").Value = "Mankaki"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 270
$ws.Range("N4").Value = 21000
$ws.Range("O4").Value = 22000
$ws.Range("P4").Value = 21500
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 1194

$ws.Range("D5").Value = 44355
$ws.Range("K5").Value = "Mankaki"
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 270
$ws.Range("N5").Value = 20000
$ws.Range("O5").Value = 21000
$ws.Range("P5").Value = 20500
$ws.Range("R5").Value = "Región Metropolitana"
$ws.Range("S5").Value = 1139

$ws.Range("D6").Value = 44301
$ws.Range("K6").Value = "Hachiya"
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 250
$ws.Range("N6").Value = 20000
$ws.Range("O6").Value = 21000
$ws.Range("P6").Value = 20500
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 1139
